$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $text)
    if ($text -match "^\s*[+-]?((\d+\.?\d*)|(\.\d+))([eE][+-]?\d+)?\s*$") {
        $origStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value2 = $text
        $cell.Style = $origStyle
    } else {
        $cell.Value2 = $text
    }
}

# --- Simple D/E value updates ---
Set-TextValue $ws.Range("D2") "22.217.87"
Set-TextValue $ws.Range("E2") "  +7.99%  "
Set-TextValue $ws.Range("D3") "1.587.15"
Set-TextValue $ws.Range("E3") "  +7.69%  "
Set-TextValue $ws.Range("D4") "1.002"
Set-TextValue $ws.Range("E4") "  -0.38%  "
Set-TextValue $ws.Range("D5") "0.9908"
Set-TextValue $ws.Range("E5") "  +3.36%  "
Set-TextValue $ws.Range("D6") "297.65"
Set-TextValue $ws.Range("E6") "  +7.44%  "
Set-TextValue $ws.Range("D8") "0.3332"
Set-TextValue $ws.Range("E8") "  +8.03%  "
Set-TextValue $ws.Range("D9") "41.16"
Set-TextValue $ws.Range("E9") "  +3.70%  "
Set-TextValue $ws.Range("D10") "1.111"
Set-TextValue $ws.Range("E10") "  +3.52%  "
Set-TextValue $ws.Range("D11") "0.06916"
Set-TextValue $ws.Range("E11") "  +3.99%  "
Set-TextValue $ws.Range("D12") "1.002"
Set-TextValue $ws.Range("E12") "  +0.12%  "
Set-TextValue $ws.Range("D13") "19.31"
Set-TextValue $ws.Range("E13") "  +5.80%  "
Set-TextValue $ws.Range("D14") "5.790"
Set-TextValue $ws.Range("E14") "  +5.16%  "
Set-TextValue $ws.Range("D15") "6.507"
Set-TextValue $ws.Range("E15") "  +5.52%  "
Set-TextValue $ws.Range("D16") "0.9917"
Set-TextValue $ws.Range("E16") "  +3.25%  "
Set-TextValue $ws.Range("D19") "0.06574"
Set-TextValue $ws.Range("E19") "  +10.67%  "
Set-TextValue $ws.Range("D20") "75.81"
Set-TextValue $ws.Range("E20") "  +9.93%  "
Set-TextValue $ws.Range("D21") "15.77"
Set-TextValue $ws.Range("E21") "  +8.38%  "
Set-TextValue $ws.Range("D22") "5.889"
Set-TextValue $ws.Range("E22") "  +6.95%  "
Set-TextValue $ws.Range("D24") "22.164.22"
Set-TextValue $ws.Range("E24") "  +7.68%  "
Set-TextValue $ws.Range("D25") "2.371"
Set-TextValue $ws.Range("E25") "  +4.80%  "
Set-TextValue $ws.Range("D26") "2.487"
Set-TextValue $ws.Range("E26") "  +16.59%  "
Set-TextValue $ws.Range("D27") "148.27"
Set-TextValue $ws.Range("E27") "  +4.23%  "
Set-TextValue $ws.Range("D28") "19.10"
Set-TextValue $ws.Range("E28") "  +11.11%  "
Set-TextValue $ws.Range("D29") "1.756.19"
Set-TextValue $ws.Range("E29") "  +7.36%  "
Set-TextValue $ws.Range("D30") "121.33"
Set-TextValue $ws.Range("E30") "  +6.54%  "
Set-TextValue $ws.Range("D31") "3.918"
Set-TextValue $ws.Range("E31") "  -0.07%  "
Set-TextValue $ws.Range("D32") "5.773"
Set-TextValue $ws.Range("E32") "  +16.65%  "
Set-TextValue $ws.Range("D33") "0.9118"
Set-TextValue $ws.Range("E33") "  +13.62%  "
Set-TextValue $ws.Range("D34") "0.08121"
Set-TextValue $ws.Range("E34") "  +1.21%  "
Set-TextValue $ws.Range("D35") "1.629"
Set-TextValue $ws.Range("E35") "  +7.97%  "
Set-TextValue $ws.Range("D37") "5.075"
Set-TextValue $ws.Range("E37") "  +7.20%  "
Set-TextValue $ws.Range("D38") "1.230"
Set-TextValue $ws.Range("E38") "  +1.31%  "
Set-TextValue $ws.Range("D39") "0.05982"
Set-TextValue $ws.Range("E39") "  +3.63%  "
Set-TextValue $ws.Range("D40") "8.256"
Set-TextValue $ws.Range("E40") "  +10.84%  "
Set-TextValue $ws.Range("D41") "0.02165"
Set-TextValue $ws.Range("E41") "  +5.49%  "
Set-TextValue $ws.Range("D44") "0.5728"
Set-TextValue $ws.Range("E44") "  +8.18%  "
Set-TextValue $ws.Range("D45") "3.749"
Set-TextValue $ws.Range("E45") "  +6.29%  "
Set-TextValue $ws.Range("D46") "12.66"
Set-TextValue $ws.Range("E46") "  +3.84%  "
Set-TextValue $ws.Range("D47") "124.48"
Set-TextValue $ws.Range("E47") "  +4.83%  "
Set-TextValue $ws.Range("D48") "0.5543"
Set-TextValue $ws.Range("E48") "  +6.29%  "
Set-TextValue $ws.Range("D49") "1.929"
Set-TextValue $ws.Range("E49") "  +6.15%  "
Set-TextValue $ws.Range("D50") "0.06708"
Set-TextValue $ws.Range("E50") "  +3.70%  "
Set-TextValue $ws.Range("D51") "72.06"
Set-TextValue $ws.Range("E51") "  +7.10%  "

# --- E-only value updates ---
Set-TextValue $ws.Range("E7") "  -0.42%  "
Set-TextValue $ws.Range("E23") "  +3.05%  "
Set-TextValue $ws.Range("E36") "  +11.81%  "

# --- Row 17/18 swap (ShibaInu <-> WrappedEther) with updated values ---
$ws.Range("B17").Value2 = "WrappedEther"
$ws.Range("C17").Value2 = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D17") "1.587.05"
Set-TextValue $ws.Range("E17") "  +7.73%  "
$ws.Range("B18").Value2 = "ShibaInu"
$ws.Range("C18").Value2 = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D18") "0.00001059"
Set-TextValue $ws.Range("E18") "  +3.10%  "

# --- Row 42/43 swap (Algorand <-> Frax) with updated values ---
$ws.Range("B42").Value2 = "Frax"
$ws.Range("C42").Value2 = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D42") "0.9912"
Set-TextValue $ws.Range("E42") "  +3.15%  "
$ws.Range("B43").Value2 = "Algorand"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D43") "0.1967"
Set-TextValue $ws.Range("E43") "  +4.56%  "
